$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

$ws.Range("H40").Value = 2909.7778
$ws.Range("I40").Value = 2802.5
$ws.Range("J40").Value = 2995.6
$ws.Range("K40").Value = 2802.5
$ws.Range("L40").Value = 2995.6
$ws.Range("M40").Value = -2627.5
$ws.Range("N40").Value = -3345.6

$ws.Range("H55").Value = 274.4
$ws.Range("I55").Value = 258
$ws.Range("J55").Value = 299
$ws.Range("K55").Value = 258
$ws.Range("L55").Value = 299
$ws.Range("M55").Value = -44
$ws.Range("N55").Value = -727

$ws.Range("H62").Value = 2317.25
$ws.Range("I62").Value = 2100
$ws.Range("J62").Value = 2969
$ws.Range("K62").Value = 2100
$ws.Range("L62").Value = 2969
$ws.Range("M62").Value = -1476
$ws.Range("N62").Value = -4217

$ws.Range("H65").Value = 2317.25
$ws.Range("I65").Value = 2100
$ws.Range("J65").Value = 2969
$ws.Range("K65").Value = 10500
$ws.Range("L65").Value = 14845
$ws.Range("M65").Value = -7380
$ws.Range("N65").Value = -21085

$ws.Range("H112").Value = 3307.7058
$ws.Range("J112").Value = 3307.7058
$ws.Range("L112").Value = 9923.117400000001
$ws.Range("N112").Value = -12139.1174

$ws.Range("H138").Value = 2551.4285
$ws.Range("I138").Value = 2973.524
$ws.Range("K138").Value = 8920.572
$ws.Range("M138").Value = -3780.572

$ws.Range("H141").Value = 3827.3572
$ws.Range("I141").Value = 2747.25
$ws.Range("J141").Value = 5267.5
$ws.Range("K141").Value = 8241.75
$ws.Range("L141").Value = 15802.5
$ws.Range("M141").Value = -3061.75
$ws.Range("N141").Value = -26162.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5577.1
$ws.Range("I32").Value = 4000.3684
$ws.Range("J32").Value = 8300.546
$ws.Range("K32").Value = 4000.3684
$ws.Range("L32").Value = 8300.546
$ws.Range("M32").Value = -3713.3684
$ws.Range("N32").Value = -8874.546

$ws.Range("H88").Value = 2830.9412
$ws.Range("I88").Value = 2153.375
$ws.Range("J88").Value = 3433.2222
$ws.Range("K88").Value = 2153.375
$ws.Range("L88").Value = 3433.2222
$ws.Range("M88").Value = -1747.375
$ws.Range("N88").Value = -4245.2222

$ws.Range("H91").Value = 2830.9412
$ws.Range("I91").Value = 2153.375
$ws.Range("J91").Value = 3433.2222
$ws.Range("K91").Value = 2153.375
$ws.Range("L91").Value = 3433.2222
$ws.Range("M91").Value = -749.375
$ws.Range("N91").Value = -6241.2222

$ws.Range("H109").Value = 67229.875
$ws.Range("J109").Value = 67229.875
$ws.Range("L109").Value = 67229.875
$ws.Range("N109").Value = -70003.875

$ws.Range("H112").Value = 24998.25
$ws.Range("J112").Value = 24998.25
$ws.Range("L112").Value = 24998.25
$ws.Range("N112").Value = -27952.25

$ws.Range("H132").Value = 2314.1
$ws.Range("I132").Value = 1791.8462
$ws.Range("J132").Value = 3284
$ws.Range("K132").Value = 5375.5386
$ws.Range("L132").Value = 9852
$ws.Range("M132").Value = -2845.5386
$ws.Range("N132").Value = -14912

$ws.Range("H135").Value = 30728.857
$ws.Range("J135").Value = 30728.857
$ws.Range("L135").Value = 30728.857
$ws.Range("N135").Value = -40868.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 99994.664
$ws.Range("J110").Value = 99994.664
$ws.Range("L110").Value = 99994.664
$ws.Range("N110").Value = -108174.664

$ws.Range("H130").Value = 59997.2
$ws.Range("J130").Value = 59997.2
$ws.Range("L130").Value = 59997.2
$ws.Range("N130").Value = -70037.2

$ws.Range("H135").Value = 34018
$ws.Range("J135").Value = 33821.6
$ws.Range("L135").Value = 33821.6
$ws.Range("N135").Value = -43961.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3673.6086
$ws.Range("I31").Value = 1183.75
$ws.Range("J31").Value = 6389.8184
$ws.Range("K31").Value = 1183.75
$ws.Range("L31").Value = 6389.8184
$ws.Range("M31").Value = -888.75
$ws.Range("N31").Value = -6979.8184

$ws.Range("H34").Value = 3673.6086
$ws.Range("I34").Value = 1183.75
$ws.Range("J34").Value = 6389.8184
$ws.Range("K34").Value = 1183.75
$ws.Range("L34").Value = 6389.8184
$ws.Range("M34").Value = -981.75
$ws.Range("N34").Value = -6793.8184

$ws.Range("H62").Value = 3832.8333
$ws.Range("I62").Value = 3874.25
$ws.Range("K62").Value = 3874.25
$ws.Range("M62").Value = -3250.25

$ws.Range("H65").Value = 3832.8333
$ws.Range("I65").Value = 3874.25
$ws.Range("K65").Value = 19371.25
$ws.Range("M65").Value = -16251.25

$ws.Range("H132").Value = 2964.6875
$ws.Range("I132").Value = 1203.5714
$ws.Range("K132").Value = 3610.7142
$ws.Range("M132").Value = -1080.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1991.1428
$ws.Range("I68").Value = 1988.6666
$ws.Range("K68").Value = 5965.9998
$ws.Range("M68").Value = -5154.9998

$ws.Range("H71").Value = 1991.1428
$ws.Range("I71").Value = 1988.6666
$ws.Range("K71").Value = 17897.9994
$ws.Range("M71").Value = -13841.9994

$ws.Range("H122").Value = 1420.875
$ws.Range("I122").Value = 1045.25
$ws.Range("K122").Value = 9407.25
$ws.Range("M122").Value = -6957.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 29900
$ws.Range("J46").Value = 29900
$ws.Range("L46").Value = 29900
$ws.Range("N46").Value = -30212

$ws.Range("H57").Value = 9800
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2002
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -10008
$ws.Range("N83").ClearContents()

$ws.Range("H122").Value = 2324.1
$ws.Range("I122").Value = 2224.3845
$ws.Range("J122").Value = 2509.2856
$ws.Range("K122").Value = 6673.1535
$ws.Range("L122").Value = 7527.8568
$ws.Range("M122").Value = -4223.1535
$ws.Range("N122").Value = -12427.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1634.9
$ws.Range("J46").Value = 1634.9
$ws.Range("L46").Value = 1634.9
$ws.Range("N46").Value = -2010.9

$ws.Range("H55").Value = 248.18182
$ws.Range("I55").Value = 216.33333
$ws.Range("J55").Value = 260.125
$ws.Range("K55").Value = 216.33333
$ws.Range("L55").Value = 260.125
$ws.Range("M55").Value = -43.33332999999999
$ws.Range("N55").Value = -606.125

$ws.Range("H61").Value = 2161.1904
$ws.Range("J61").Value = 2930.8
$ws.Range("L61").Value = 2930.8
$ws.Range("N61").Value = -3334.8

$ws.Range("H82").Value = 3908.7778
$ws.Range("I82").Value = 1824.5
$ws.Range("J82").Value = 4504.2856
$ws.Range("K82").Value = 1824.5
$ws.Range("L82").Value = 4504.2856
$ws.Range("M82").Value = -1463.5
$ws.Range("N82").Value = -5226.2856

$ws.Range("H85").Value = 3908.7778
$ws.Range("I85").Value = 1824.5
$ws.Range("J85").Value = 4504.2856
$ws.Range("K85").Value = 1824.5
$ws.Range("L85").Value = 4504.2856
$ws.Range("M85").Value = -576.5
$ws.Range("N85").Value = -7000.2856

$ws.Range("H104").Value = 10090.75
$ws.Range("J104").Value = 10090.75
$ws.Range("L104").Value = 10090.75
$ws.Range("N104").Value = -17078.75

$ws.Range("H110").Value = 22525
$ws.Range("J110").Value = 22525
$ws.Range("L110").Value = 22525
$ws.Range("N110").Value = -30705

$ws.Range("H113").Value = 2161.1904
$ws.Range("J113").Value = 2930.8
$ws.Range("L113").Value = 2930.8
$ws.Range("N113").Value = -7270.8

$ws.Range("H127").Value = 49715
$ws.Range("J127").Value = 49715
$ws.Range("L127").Value = 49715
$ws.Range("N127").Value = -59635

$ws.Range("H133").Value = 87326
$ws.Range("J133").Value = 87326
$ws.Range("L133").Value = 87326
$ws.Range("N133").Value = -92386

$ws.Range("H136").Value = 4406.731
$ws.Range("J136").Value = 5318.7334
$ws.Range("L136").Value = 15956.2002
$ws.Range("N136").Value = -21056.2002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 2634
$ws.Range("J14").Value = 3117.5
$ws.Range("L14").Value = 3117.5
$ws.Range("N14").Value = -3453.5

$ws.Range("H107").Value = 627.73334
$ws.Range("I107").Value = 529.7143
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1589.1429
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = 330.8571000000002
$ws.Range("N107").Value = -9840

$ws.Range("H132").Value = 1797.3125
$ws.Range("I132").Value = 978.0909
$ws.Range("K132").Value = 2934.2727
$ws.Range("M132").Value = -404.2727

$ws.Range("H136").Value = 2537.8157
$ws.Range("I136").Value = 2068.7144
$ws.Range("J136").Value = 3851.3
$ws.Range("K136").Value = 6206.1432
$ws.Range("L136").Value = 11553.9
$ws.Range("M136").Value = -3656.1432
$ws.Range("N136").Value = -16653.9
